# Import newly received loss-of-sale records into the Walk-In Report sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Each entry: #, Date, Customer Name, Contact, Function Date, Staff, Status,
#             Category, Sub Category, Repeat count, Remarks
$records = @(
    @(13, "22-12-2025", "Abhishek", 8078946799, "06-01-2026", "VISHNU N", "Loss", "ENQUIRY",            "ENQUIRY WITHOUT BRIDE/FAMILY", "-", "just checking"),
    @(14, "23-12-2025", "HASHIM",   7736621379, "27-12-2025", "VISHNU N", "Loss", "SIZE NOT SUITABLE",  "SIZE TOO SMALL",                "-", "jst cheaking conform in evening"),
    @(15, "23-12-2025", "NISHAD",   8289924057, "31-12-2025", "ARJUN P",  "Loss", "PRODUCT",             "REQUIRED MODEL NOT AVAILABLE",  "-", "kurtha"),
    @(16, "23-12-2025", "Amal",     8089143196, "26-12-2025", "VISHNU N", "Loss", "ENQUIRY",             "ENQUIRY WITHOUT TRIAL",         "-", "conform later January function"),
    @(17, "24-12-2025", "harshan",  9048520130, "25-01-2026", "ARJUN P",  "Loss", "ENQUIRY",             "Enquiry for Relative/Friend",   "-", "just checking all nearby stores"),
    @(18, "25-12-2025", "Abhiram",  9544484863, "04-01-2026", "ARJUN P",  "Loss", "ENQUIRY",             "Enquiry for Relative/Friend",   "-", "just checking and after coming")
)

$startRow = 15
for ($i = 0; $i -lt $records.Count; $i++) {
    $row = $startRow + $i
    $rec = $records[$i]

    $ws.Cells.Item($row, 1).Value2 = $rec[0]
    $ws.Cells.Item($row, 1).NumberFormat = "0"

    # Date-like text ("DD-MM-YYYY") must stay plain text, not get
    # auto-converted into a date serial by Excel's input parser. Force the
    # cell to Text format while assigning, then drop back to the sheet's
    # normal (unstyled) look so the saved cell matches the plain data rows.
    $ws.Cells.Item($row, 2).NumberFormat = "@"
    $ws.Cells.Item($row, 2).Value2 = $rec[1]
    $ws.Cells.Item($row, 2).Style = "Normal"

    $ws.Cells.Item($row, 3).Value2 = $rec[2]

    $ws.Cells.Item($row, 4).Value2 = $rec[3]
    $ws.Cells.Item($row, 4).NumberFormat = "0"

    $ws.Cells.Item($row, 5).NumberFormat = "@"
    $ws.Cells.Item($row, 5).Value2 = $rec[4]
    $ws.Cells.Item($row, 5).Style = "Normal"

    $ws.Cells.Item($row, 6).Value2 = $rec[5]
    $ws.Cells.Item($row, 7).Value2 = $rec[6]
    $ws.Cells.Item($row, 8).Value2 = $rec[7]
    $ws.Cells.Item($row, 9).Value2 = $rec[8]
    $ws.Cells.Item($row, 10).Value2 = $rec[9]
    $ws.Cells.Item($row, 11).Value2 = $rec[10]
}
